# Applies the "AdminLog Completed, FrontWeb User some" edit to the DB Design
# workbook: updates the TFit_Users table block (rows 64-73 on Sheet1) by
# removing the Age/Gender fields and appending two new fields
# (OperateCode int, IsActivated bit), and updates the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the TFit_Users field list (rows 64-73) ---
# Before:
#   64 Age                     int
#   65 Gender                  int
#   66 PhoneNum                varchar(20)
#   67 PasswordSalt            varchar(20)
#   68 PasswordHash            varchar(100)
#   69 Email                   varchar(30)
#   70 LoginErrorTimes         int
#   71 LastLoginErrorDateTime  datetime
#   72 IsDeleted               bit
#   73 CreateDateTime          datetime
# After (Age/Gender removed, everything shifts up 2, two new rows appended):
#   64 PhoneNum                varchar(20)
#   65 PasswordSalt            varchar(20)
#   66 PasswordHash            varchar(100)
#   67 Email                   varchar(30)
#   68 LoginErrorTimes         int
#   69 LastLoginErrorDateTime  datetime
#   70 IsDeleted               bit
#   71 CreateDateTime          datetime
#   72 OperateCode             int
#   73 IsActivated             bit

$ws.Range("A64").Value = "PhoneNum"
$ws.Range("B64").Value = "varchar(20)"

$ws.Range("A65").Value = "PasswordSalt"
$ws.Range("B65").Value = "varchar(20)"

$ws.Range("A66").Value = "PasswordHash"
$ws.Range("B66").Value = "varchar(100)"

$ws.Range("A67").Value = "Email"
$ws.Range("B67").Value = "varchar(30)"

$ws.Range("A68").Value = "LoginErrorTimes"
$ws.Range("B68").Value = "int"

$ws.Range("A69").Value = "LastLoginErrorDateTime"
$ws.Range("B69").Value = "datetime"

$ws.Range("A70").Value = "IsDeleted"
$ws.Range("B70").Value = "bit"

$ws.Range("A71").Value = "CreateDateTime"
$ws.Range("B71").Value = "datetime"

$ws.Range("A72").Value = "OperateCode"
$ws.Range("B72").Value = "int"

$ws.Range("A73").Value = "IsActivated"
$ws.Range("B73").Value = "bit"

# --- Update the saved sheet view / selection (scrolled to A37, cell A72 selected) ---
$ws.Range("A72").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 37
$window.ScrollColumn = 1
